$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.981.42"
$ws.Range("E2").Value = "  +2.79%  "

$ws.Range("D3").Value = "3.045.22"
$ws.Range("E3").Value = "  +2.24%  "

$ws.Range("E4").Value = "  -0.03%  "

$r = $ws.Range("D5")
$r.Value = "'519.18"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +3.24%  "

$r = $ws.Range("D6")
$r.Value = "'141.49"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +5.01%  "

$r = $ws.Range("D7")
$r.Value = "'0.999"
$r.Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "

$r = $ws.Range("D8")
$r.Value = "'0.445"
$r.Style = "Normal"

$r = $ws.Range("D9")
$r.Value = "'7.49"
$r.Style = "Normal"
$ws.Range("E9").Value = "  +2.16%  "

$r = $ws.Range("D10")
$r.Value = "'0.112"
$r.Style = "Normal"
$ws.Range("E10").Value = "  +5.87%  "

$r = $ws.Range("D11")
$r.Value = "'0.369"
$r.Style = "Normal"
$ws.Range("E11").Value = "  +5.42%  "

$ws.Range("D12").Value = "3.576.40"
$ws.Range("E12").Value = "  +2.48%  "

$ws.Range("E13").Value = "  +2.25%  "

$r = $ws.Range("D14")
$r.Value = "'26.74"
$r.Style = "Normal"
$ws.Range("E14").Value = "  +6.45%  "

$ws.Range("E15").Value = "  +13.37%  "

$ws.Range("D16").Value = "57.971.29"
$ws.Range("E16").Value = "  +2.75%  "

$r = $ws.Range("D17")
$r.Value = "'6.20"
$r.Style = "Normal"
$ws.Range("E17").Value = "  +9.59%  "

$ws.Range("D18").Value = "3.046.96"
$ws.Range("E18").Value = "  +2.30%  "

$r = $ws.Range("D19")
$r.Value = "'13.03"
$r.Style = "Normal"
$ws.Range("E19").Value = "  +5.68%  "

$r = $ws.Range("D20")
$r.Value = "'8.10"
$r.Style = "Normal"
$ws.Range("E20").Value = "  +4.33%  "

$r = $ws.Range("D21")
$r.Value = "'336.67"
$r.Style = "Normal"
$ws.Range("E21").Value = "  +3.88%  "

$r = $ws.Range("D22")
$r.Value = "'5.77"
$r.Style = "Normal"
$ws.Range("E22").Value = "  +0.96%  "

$r = $ws.Range("D23")
$r.Value = "'0.999"
$r.Style = "Normal"
$ws.Range("E23").Value = "  -0.14%  "

$r = $ws.Range("D24")
$r.Value = "'0.501"
$r.Style = "Normal"
$ws.Range("E24").Value = "  +6.70%  "

$r = $ws.Range("D25")
$r.Value = "'65.05"
$r.Style = "Normal"
$ws.Range("E25").Value = "  +5.09%  "

$r = $ws.Range("D26")
$r.Value = "'0.169"
$r.Style = "Normal"
$ws.Range("E26").Value = "  +3.80%  "

$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$r = $ws.Range("D27")
$r.Value = "'1.02"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +1.63%  "

$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0955"
$ws.Range("E28").Value = "  +7.19%  "

$r = $ws.Range("D29")
$r.Value = "'6.89"
$r.Style = "Normal"
$ws.Range("E29").Value = "  +5.78%  "

$r = $ws.Range("D30")
$r.Value = "'7.52"
$r.Style = "Normal"
$ws.Range("E30").Value = "  +11.20%  "

$ws.Range("E31").Value = "  +5.01%  "

$r = $ws.Range("D32")
$r.Value = "'1.23"
$r.Style = "Normal"
$ws.Range("E32").Value = "  +2.71%  "

$r = $ws.Range("D33")
$r.Value = "'21.06"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +3.18%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$r = $ws.Range("D34")
$r.Value = "'4.76"
$r.Style = "Normal"
$ws.Range("E34").Value = "  +7.12%  "

$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$r = $ws.Range("D35")
$r.Value = "'156.52"
$r.Style = "Normal"
$ws.Range("E35").Value = "  -0.18%  "

$ws.Range("E36").Value = "  +7.15%  "

$ws.Range("E37").Value = "  +1.78%  "

$ws.Range("E38").Value = "  +9.52%  "

$r = $ws.Range("D39")
$r.Value = "'0.0691"
$r.Style = "Normal"
$ws.Range("E39").Value = "  +2.83%  "

$ws.Range("D40").Value = "3.081.67"
$ws.Range("E40").Value = "  +2.24%  "

$r = $ws.Range("D41")
$r.Value = "'37.67"
$r.Style = "Normal"
$ws.Range("E41").Value = "  +4.08%  "

$r = $ws.Range("D42")
$r.Value = "'3.90"
$r.Style = "Normal"
$ws.Range("E42").Value = "  +9.73%  "

$ws.Range("E43").Value = "  +0.09%  "

$ws.Range("E44").Value = "  +3.88%  "

$ws.Range("D45").Value = "2.315.88"
$ws.Range("E45").Value = "  +3.02%  "

$ws.Range("E46").Value = "  +4.21%  "

$r = $ws.Range("D47")
$r.Value = "'1.01"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +2.33%  "

$r = $ws.Range("D48")
$r.Value = "'6.07"
$r.Style = "Normal"
$ws.Range("E48").Value = "  +5.83%  "

$ws.Range("E49").Value = "  +2.45%  "

$r = $ws.Range("D50")
$r.Value = "'19.73"
$r.Style = "Normal"
$ws.Range("E50").Value = "  +4.16%  "

$ws.Range("E51").Value = "  -4.41%  "
